$d = $word.ActiveDocument

# Locate the paragraph containing the "LOQ4209" requirement text.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*LOQ4209: Engenharia da Qualidade I (Requisito fraco)*") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    # Delete the three paragraphs that follow it:
    #   - an empty paragraph
    #   - "Ver no Jupiter Salvar em pdf Salvar em docx"
    #   - the "© 2020 . Contact: ..." paragraph
    $pStart = $d.Paragraphs.Item($target + 1)
    $pEnd = $d.Paragraphs.Item($target + 3)

    $r = $d.Range($pStart.Range.Start, $pEnd.Range.End)
    $r.Delete()
}
